# Generate Report for Handback
# Refresh the timestamps recorded on the handback status report:
#  - Overview!G2      "Latest HO Xliff Generate Date"
#  - zh-cn!H2          "Correspond Handoff Datetime"
#  - zh-cn!K2          "Correspond Handback DateTime"
#  - de-de!K2          "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 13:01:49"
$wsOverview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 13:01:44"
$wsZhCn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K2").Value = "2016-08-27 13:02:04"
$wsZhCn.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-27 13:02:12"
$wsDeDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
